$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New rows of resource data (A2 + B2:B22) ---------------------------
# A2 holds the script code, B2:B22 holds the resource numbers REC10..REC210
# (a simple +10 series), mirroring a fill-down/typed list in the sheet.
$ws.Range("A2").Value = "CS_08_05_CO"

$resourceNumbers = @(
    "REC10","REC20","REC30","REC40","REC50","REC60","REC70","REC80","REC90",
    "REC100","REC110","REC120","REC130","REC140","REC150","REC160","REC170",
    "REC180","REC190","REC200","REC210"
)
for ($i = 0; $i -lt $resourceNumbers.Length; $i++) {
    $ws.Cells.Item($i + 2, 2).Value = $resourceNumbers[$i]
}

# --- Column width tweaks to make room for the new columns ---------------
# (B and D are brand-new data columns; H and I were widened by hand.)
$ws.Columns.Item(2).ColumnWidth = 17.42578125 - 0.8333333333333334
$ws.Columns.Item(4).ColumnWidth = 13.7109375 - 0.8333333333333334
$ws.Columns.Item(8).ColumnWidth = 24.42578125 - 0.8333333333333334
$ws.Columns.Item(9).ColumnWidth = 22.28515625 - 0.8333333333333334

# --- Selection left where the author's cursor ended up ------------------
$ws.Range("C2").Select()
